$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: drop trailing significance star (one star removed from each)
$ws.Range("B2").Value = "0.011`n (0.019)"
$ws.Range("C2").Value = "0.076`n (0.062)"
$ws.Range("D2").Value = "0.007`n (0.012)"
$ws.Range("E2").Value = "0.012`n (0.037)"
$ws.Range("F2").Value = "-0.017*`n (0.010)"
$ws.Range("G2").Value = "-0.006`n (0.026)"

# Row 3
$ws.Range("B3").Value = "-0.391***`n (0.102)"
$ws.Range("D3").Value = "0.516***`n (0.067)"
$ws.Range("F3").Value = "0.380***`n (0.052)"
$ws.Range("G3").Value = "0.102`n (0.090)"

# Row 4
$ws.Range("B4").Value = "0.291*`n (0.174)"
$ws.Range("D4").Value = "0.071`n (0.113)"
$ws.Range("F4").Value = "0.168*`n (0.094)"
